$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 66: GenomeWeb link for the new FDA companion-diagnostic reclassification story ---
$ws.Range("A66").Value = "https://www.genomeweb.com/cancer/fda-proposes-reclassification-companion-diagnostic-tests"
$ws.Hyperlinks.Add($ws.Range("A66"), "https://www.genomeweb.com/cancer/fda-proposes-reclassification-companion-diagnostic-tests")
$ws.Range("A65").Copy()
$ws.Range("A66").PasteSpecial(-4122)
$ws.Range("B66").Value = "companion diagnostic"
$ws.Range("C66").Value = "FDA Proposes Reclassification of Companion Diagnostic Tests"

# --- Row 67: matching 360dx link for the same story ---
$ws.Range("A67").Value = "https://www.360dx.com/cancer/fda-proposes-reclassification-companion-diagnostic-tests"
$ws.Hyperlinks.Add($ws.Range("A67"), "https://www.360dx.com/cancer/fda-proposes-reclassification-companion-diagnostic-tests")
$ws.Range("A65").Copy()
$ws.Range("A67").PasteSpecial(-4122)
$ws.Range("B67").Value = "companion diagnostic"
$ws.Range("C67").Value = "FDA Proposes Reclassification of Companion Diagnostic Tests"

$excel.CutCopyMode = 0
